$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# ALC sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

# Row 134
$ws.Range("H134").Value = 65274.668
$ws.Range("I134").Value = 20000
$ws.Range("K134").Value = 20000
$ws.Range("M134").Value = -14930

# Row 138
$ws.Range("H138").Value = 4437.0806
$ws.Range("I138").Value = 2570.0625
$ws.Range("J138").Value = 4796.988
$ws.Range("K138").Value = 7710.1875
$ws.Range("L138").Value = 14390.964
$ws.Range("M138").Value = -2570.1875
$ws.Range("N138").Value = -24670.964

# ---------------------------------------------------------------------------
# ARM sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

# Row 121
$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 0

# Row 122
$ws.Range("H122").Value = 1100
$ws.Range("I122").Value = 1032.3077
$ws.Range("J122").Value = 1276
$ws.Range("K122").Value = 3096.9231
$ws.Range("L122").Value = 3828
$ws.Range("M122").Value = -646.9231
$ws.Range("N122").Value = -8728

# Row 123
$ws.Range("H123").Value = 39700
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 39700
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 39700
$ws.Range("N123").Value = -49500

# Row 124
$ws.Range("H124").Value = 43643
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 43643
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 43643
$ws.Range("N124").Value = -53463

# Row 125
$ws.Range("H125").Value = 36137.145
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 36137.145
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 36137.145
$ws.Range("N125").Value = -45977.145

# Row 126
$ws.Range("H126").Value = 7420
$ws.Range("I126").Value = 7420
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 22260
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -19790

# Row 127
$ws.Range("H127").Value = 50993.332
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 50993.332
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 50993.332
$ws.Range("N127").Value = -60913.332

# Row 128
$ws.Range("H128").Value = 49980
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 49980
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 49980
$ws.Range("N128").Value = -59940

# Row 129
$ws.Range("H129").Value = 49999
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 49999
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 49999
$ws.Range("N129").Value = -59999

# Row 130
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0

# Row 131
$ws.Range("H131").Value = 59980
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 59980
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 59980
$ws.Range("N131").Value = -70060

# Row 132
$ws.Range("H132").Value = 2690.5952
$ws.Range("I132").Value = 2442.8708
$ws.Range("J132").Value = 3388.7273
$ws.Range("K132").Value = 7328.6124
$ws.Range("L132").Value = 10166.1819
$ws.Range("M132").Value = -4798.6124
$ws.Range("N132").Value = -15226.1819

# Row 133
$ws.Range("H133").Value = 35000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 35000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 35000
$ws.Range("N133").Value = -40060

# Row 134
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0

# Row 135
$ws.Range("H135").Value = 35346.855
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 35346.855
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 35346.855
$ws.Range("N135").Value = -45486.855

# Row 137
$ws.Range("H137").Value = 73400
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 73400
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 73400
$ws.Range("N137").Value = -83600

# Row 138
$ws.Range("H138").Value = 50000
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 50000
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60280

# Row 139
$ws.Range("H139").Value = 28333.5
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 28333.5
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 28333.5
$ws.Range("N139").Value = -38613.5

# Row 140
$ws.Range("H140").Value = 54232.25
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 54232.25
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 54232.25
$ws.Range("N140").Value = -64592.25

# Row 141
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0

# ---------------------------------------------------------------------------
# CRP sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

# Row 4
$ws.Range("H4").Value = 6900
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 6900
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 6900
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -7124

# ---------------------------------------------------------------------------
# CUL sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

# Row 131
$ws.Range("H131").Value = 14737455
$ws.Range("I131").Value = 55667036
$ws.Range("J131").Value = 2806.32
$ws.Range("K131").Value = 167001108
$ws.Range("L131").Value = 8418.960000000001
$ws.Range("M131").Value = -166996068
$ws.Range("N131").Value = -18498.96

# ---------------------------------------------------------------------------
# LTW sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

# Row 82
$ws.Range("H82").Value = 1883.3889
$ws.Range("I82").Value = 1607.8462
$ws.Range("J82").Value = 2599.8
$ws.Range("K82").Value = 1607.8462
$ws.Range("L82").Value = 2599.8
$ws.Range("M82").Value = -1246.8462
$ws.Range("N82").Value = -3321.8

# Row 85
$ws.Range("H85").Value = 1883.3889
$ws.Range("I85").Value = 1607.8462
$ws.Range("J85").Value = 2599.8
$ws.Range("K85").Value = 1607.8462
$ws.Range("L85").Value = 2599.8
$ws.Range("M85").Value = -359.8462
$ws.Range("N85").Value = -5095.8

# ---------------------------------------------------------------------------
# WVR sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

# Row 136
$ws.Range("H136").Value = 3319.853
$ws.Range("I136").Value = 4653
$ws.Range("J136").Value = 2267.3684
$ws.Range("K136").Value = 13959
$ws.Range("L136").Value = 6802.1052
$ws.Range("M136").Value = -11409
$ws.Range("N136").Value = -11902.1052
